$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.245.94"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3
$ws.Range("D3").Value = "3.341.98"
$ws.Range("E3").Value = "  +0.39%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'585.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.53%  "

# Row 6
$ws.Range("E6").Value = "  -1.78%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").Value = "3.338.65"
$ws.Range("E8").Value = "  +0.49%  "

# Row 9
$ws.Range("E9").Value = "  -2.51%  "

# Row 10
$ws.Range("E10").Value = "  -2.12%  "

# Row 11
$ws.Range("E11").Value = "  -1.30%  "

# Row 12
$ws.Range("D12").Value = "'46.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.83%  "

# Row 13
$ws.Range("E13").Value = "  -1.48%  "

# Row 14
$ws.Range("D14").Value = "'668.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.11%  "

# Row 15
$ws.Range("D15").Value = "3.876.39"
$ws.Range("E15").Value = "  +0.48%  "

# Row 16
$ws.Range("D16").Value = "'8.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.24%  "

# Row 17
$ws.Range("D17").Value = "66.400.50"
$ws.Range("E17").Value = "  +0.18%  "

# Row 18
$ws.Range("E18").Value = "  -0.62%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.342.57"
$ws.Range("E19").Value = "  +0.06%  "

# Row 20
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'17.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "

# Row 21
$ws.Range("D21").Value = "'11.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "

# Row 22
$ws.Range("D22").Value = "'0.897"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.16%  "

# Row 23
$ws.Range("D23").Value = "'17.68"
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'101.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.49%  "

# Row 25
$ws.Range("D25").Value = "'5.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.23%  "

# Row 26
$ws.Range("E26").Value = "  -0.46%  "

# Row 27
$ws.Range("E27").Value = "  -0.25%  "

# Row 28
$ws.Range("D28").Value = "'9.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.92%  "

# Row 29
$ws.Range("D29").Value = "'32.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.49%  "

# Row 30
$ws.Range("E30").Value = "  -2.17%  "

# Row 31
$ws.Range("D31").Value = "'6.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.41%  "

# Row 32
$ws.Range("D32").Value = "'615.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.35%  "

# Row 33
$ws.Range("E33").Value = "  +0.70%  "

# Row 34
$ws.Range("D34").Value = "'11.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.66%  "

# Row 35
$ws.Range("D35").Value = "3.859.26"
$ws.Range("E35").Value = "  +3.47%  "

# Row 36
$ws.Range("E36").Value = "  -0.81%  "

# Row 37
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
$ws.Range("D38").Value = "'56.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.88%  "

# Row 39
$ws.Range("E39").Value = "  -2.83%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0701"
$ws.Range("E40").Value = "  -4.28%  "

# Row 41
$ws.Range("E41").Value = "  -1.19%  "

# Row 42
$ws.Range("D42").Value = "'32.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.90%  "

# Row 43
$ws.Range("E43").Value = "  -3.25%  "

# Row 44
$ws.Range("D44").Value = "'3.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.70%  "

# Row 45
$ws.Range("E45").Value = "  -2.50%  "

# Row 46
$ws.Range("E46").Value = "  -2.03%  "

# Row 47
$ws.Range("D47").Value = "'2.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -15.27%  "

# Row 48
$ws.Range("E48").Value = "  -1.84%  "

# Row 49
$ws.Range("E49").Value = "  +0.29%  "

# Row 50
$ws.Range("E50").Value = "  -2.31%  "

# Row 51
$ws.Range("D51").Value = "'1.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.13%  "

